$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 54.78088418710104
$ws.Range("B3").Value = 0.9306427053331269
$ws.Range("B4").Value = 0.05572639886396525
$ws.Range("B5").Value = 0.4148140362643069
